$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3 value (Hydrogen / Iron & steel)
$ws.Range("B3").Value = 3340334.517681919

# Clear D3 (Hydrogen / Non-metallic minerals) - becomes blank
$ws.Range("D3").Value = ""

# Update C4 (Methanol / Chemicals)
$ws.Range("C4").Value = 38.9303125892257

# Update C5 (Ammonia / Chemicals)
$ws.Range("C5").Value = 0

# Row 7: rename "Other" -> "Biogas", update D7
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 342.2132282291858

# New Row 8: "Other" with D8 value; copy style from A7 to A8
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 1750.928536564099
